# Versión para Abril - Diciembre 2024
# Adjust vacation days, assigned load and computed monthly ratios per the
# commit's manual re-assignment of the schedule for Apr-Dec 2024.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Gomez
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 6.224379325864565
$ws.Range("G2").Value = 20
$ws.Range("O2").Value = 0.2512396694214876

# Row 3 - Bravo
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 6.224379325864565

# Row 4 - Iñiguez
$ws.Range("E4").Value = 0

# Row 5 - Breinbauer
$ws.Range("E5").Value = 0

# Row 6 - Arredondo
$ws.Range("E6").Value = 0

# Row 7 - Carrasco
$ws.Range("E7").Value = 0

# Row 8 - Culaciati
$ws.Range("E8").Value = 0

# Row 9 - Contreras
$ws.Range("E9").Value = 0

# Row 10 - Cisternas
$ws.Range("E10").Value = 0
$ws.Range("G10").Value = 30
$ws.Range("P10").Value = 0.2512396694214876

# Row 11 - Pio
$ws.Range("E11").Value = 0

# Row 12 - Alvo
$ws.Range("E12").Value = 0

# Row 13 - Boettiger
$ws.Range("E13").Value = 0
$ws.Range("G13").Value = 50
$ws.Range("R13").Value = 0.2512396694214876

# Row 14 - Loch
$ws.Range("E14").Value = 0

# Row 15 - Rubio
$ws.Range("E15").Value = 0
